# Updates the cryptos list (price + 1h volume change) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.596.54"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.422.33"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.01"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.06"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.46"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.795.80"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.406.77"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.632.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.69"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.44"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.29"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("E31").Value = "  +17.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.87%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +3.02%  "
$ws.Range("E36").Value = "  +2.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "130.20"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +27.83%  "
$ws.Range("E38").Value = "  +4.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.37"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.108"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.92"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.944.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("E46").Value = "  +2.68%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.658.58"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("E49").Value = "  +3.94%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  +0.10%  "
